$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9995338795734406
$ws.Range("E2").Value = 0.9995338795734406

$ws.Range("D3").Value = 0.9998343273825507
$ws.Range("E3").Value = 0.9998343273825507

$ws.Range("D4").Value = 0.00001167657564129454
$ws.Range("E4").Value = 0.00001167657564129454

$ws.Range("D5").Value = 0.137056723722291
$ws.Range("E5").Value = 0.137056723722291

$ws.Range("D6").Value = 0.9164526429078863
$ws.Range("E6").Value = 0.9164526429078863

$ws.Range("D7").Value = 0.9999971045984611
$ws.Range("E7").Value = 0.000002895401538860476

$ws.Range("D8").Value = 0.9999870893260018
$ws.Range("E8").Value = 0.00001291067399822321

$ws.Range("D9").Value = 0.9999999998099465
$ws.Range("E9").Value = 0.0000000001900535284704574

$ws.Range("D10").Value = 0.9053378290234583
$ws.Range("E10").Value = 0.09466217097654173

$ws.Range("D11").Value = 0.999998564717348
$ws.Range("E11").Value = 0.000001435282651951475
$ws.Range("F11").Value = 1.910578966140747
